$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: insert a new record at row 55 (pushing existing
# rows 55-148 down to 56-149) and populate it with this week's data point.
$ws.Rows.Item(55).Insert()

$ws.Cells.Item(55, 1).Value = 1
$ws.Cells.Item(55, 2).Value = 'Agrícola del Norte S.A. de Arica'
$ws.Cells.Item(55, 3).Value = 'Arica y Parinacota'
$ws.Cells.Item(55, 4).Value = Get-Date -Year 2021 -Month 10 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(55, 5).Value = 15
$ws.Cells.Item(55, 6).Value = 'Fruta'
$ws.Cells.Item(55, 7).Value = 100108
$ws.Cells.Item(55, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(55, 9).Value = 100108006
$ws.Cells.Item(55, 10).Value = 'Plátano'
$ws.Cells.Item(55, 11).Value = 'Sin especificar'
$ws.Cells.Item(55, 12).Value = 'Pintón'
$ws.Cells.Item(55, 13).Value = 120
$ws.Cells.Item(55, 14).Value = 22000
$ws.Cells.Item(55, 15).Value = 23000
$ws.Cells.Item(55, 16).Value = 22500
$ws.Cells.Item(55, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(55, 18).Value = 'Bolivia'
$ws.Cells.Item(55, 19).Value = 1125
$ws.Cells.Item(55, 20).Value = 20
